$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Teste"
$ws.Range("A2").Value = "Excel"
$ws.Range("A3").Value = "para"
$ws.Range("A4").Value = "Desafio de Projeto"
$ws.Range("A5").Value = "DIO"

$null = $ws.Range("A6").Select()
